# Weekly update: insert a new weekly price record (row) for Orégano at
# "Vega Central Mapocho de Santiago" ahead of the existing history, shifting
# the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 63; rows 63-68 shift down to 64-69.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the latest weekly record.
$ws.Cells.Item(63, 1).Value = 9
$ws.Cells.Item(63, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44748
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 100112029
$ws.Cells.Item(63, 7).Value = "Orégano"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 16
$ws.Cells.Item(63, 11).Value = 16000
$ws.Cells.Item(63, 12).Value = 16000
$ws.Cells.Item(63, 13).Value = 16000
$ws.Cells.Item(63, 14).Value = "$/docena de atados"
$ws.Cells.Item(63, 15).Value = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value = 5333
$ws.Cells.Item(63, 17).Value = 3
$ws.Cells.Item(63, 18).Value = "Hortaliza"
